# "cambios de local a nube" - add the Spanish sentence
# "Hola estoy trabajando localmente" into the (until now empty) first
# paragraph, reproducing the run/proofErr structure Word itself would
# generate when the sentence is typed with the English (en-US) proofing
# language active (estoy / trabajando flagged as misspellings).

$d = $word.ActiveDocument
$p = $d.Paragraphs(1)
$r = $p.Range
$r.Collapse(1)   # wdCollapseStart - insert at the start of the paragraph

$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512">
    <pkg:xmlData>
      <Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships">
        <Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/>
      </Relationships>
    </pkg:xmlData>
  </pkg:part>
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r>
              <w:rPr>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t xml:space="preserve">Hola </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:rPr>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t>estoy</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:rPr>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:rPr>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t>trabajando</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:rPr>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t xml:space="preserve"> localmente</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$r.InsertXML($xml)
